$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "43.508.09"
$ws.Cells.Item(2, 5).Value = "  +0.26%  "

$ws.Cells.Item(3, 4).Value = "2.335.04"
$ws.Cells.Item(3, 5).Value = "  -1.67%  "

$ws.Cells.Item(4, 5).Value = "  +0.01%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "304.36"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -1.80%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "101.34"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -3.45%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.512"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -2.01%  "

$ws.Cells.Item(8, 5).Value = "  +0.05%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.514"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -0.95%  "

$ws.Cells.Item(10, 5).Value = "  -3.02%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "51.62"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -3.40%  "

$ws.Cells.Item(12, 5).Value = "  -2.14%  "

$ws.Cells.Item(13, 5).Value = "  +0.28%  "

$ws.Cells.Item(15, 4).Value = "2.701.40"
$ws.Cells.Item(15, 5).Value = "  -1.45%  "

$ws.Cells.Item(16, 5).Value = "  -0.17%  "

$ws.Cells.Item(17, 4).Value = "2.370.84"
$ws.Cells.Item(17, 5).Value = "  -0.04%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.806"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -1.40%  "

$ws.Cells.Item(19, 4).Value = "43.403.27"
$ws.Cells.Item(19, 5).Value = "  +0.08%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "11.76"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -2.10%  "

$ws.Cells.Item(21, 5).Value = "  -1.70%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "6.11"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -2.42%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "67.89"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -0.76%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "238.39"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -1.51%  "

$ws.Cells.Item(25, 5).Value = "  -3.83%  "

$ws.Cells.Item(26, 5).Value = "  -3.36%  "

$ws.Cells.Item(27, 5).Value = "  -0.19%  "

$ws.Cells.Item(28, 5).Value = "  -3.61%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "34.60"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -6.65%  "

$ws.Cells.Item(30, 5).Value = "  -2.26%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "9.24"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -3.81%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "164.77"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +1.69%  "

$ws.Cells.Item(33, 5).Value = "  -0.06%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "5.06"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -4.18%  "

$ws.Cells.Item(35, 5).Value = "  -4.93%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "4.51"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -5.44%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "16.89"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -7.69%  "

$ws.Cells.Item(38, 5).Value = "  -4.88%  "

$ws.Cells.Item(39, 5).Value = "  -7.25%  "

$ws.Cells.Item(40, 5).Value = "  -6.42%  "

$ws.Cells.Item(41, 5).Value = "  -3.17%  "

$ws.Cells.Item(42, 5).Value = "  -3.03%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "2.40"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -3.38%  "

$ws.Cells.Item(44, 4).Value = "1.980.75"
$ws.Cells.Item(44, 5).Value = "  -1.25%  "

$ws.Cells.Item(45, 5).Value = "  -1.84%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "18.60"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -6.92%  "

$ws.Cells.Item(47, 5).Value = "  -7.07%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "9.92"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -5.56%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "55.88"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -4.20%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "4.92"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +4.45%  "

$ws.Cells.Item(51, 4).Value = "2.561.83"
$ws.Cells.Item(51, 5).Value = "  +0.27%  "
